$wb = $excel.ActiveWorkbook

# --- Repayment Schedule sheet: insert a new (blank) column before column N ---
# This shifts the old "Late" column (N) to O, and the old "Outstanding" values
# column (P) to Q, leaving a new blank column N in between (RBI / Variable
# Instalments layout change).
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns.Item(14).Insert()

# Restore the width that Excel carries over to the newly inserted column
# (matches the width the "In Advance" column to its left used to have).
$ws.Columns.Item(14).ColumnWidth = 10.25

# --- Make "Repayment Schedule" the active sheet/tab, with a new selection ---
$null = $ws.Activate()
$ws.Range("L19").Select() | Out-Null
